$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "84+7="
$t.Cell(1,2).Range.Text = "73-70="
$t.Cell(1,3).Range.Text = "90-84="
$t.Cell(1,4).Range.Text = "84-20="
$t.Cell(1,5).Range.Text = "58-44="

$t.Cell(2,1).Range.Text = "52-29="
$t.Cell(2,2).Range.Text = "57-50="
$t.Cell(2,3).Range.Text = "14+13="
$t.Cell(2,4).Range.Text = "40+1="
$t.Cell(2,5).Range.Text = "17+59="

$t.Cell(3,1).Range.Text = "26+60="
$t.Cell(3,2).Range.Text = "33-14="
$t.Cell(3,3).Range.Text = "26+33="
$t.Cell(3,4).Range.Text = "88-60="
$t.Cell(3,5).Range.Text = "61+5="

$t.Cell(4,1).Range.Text = "79-24="
$t.Cell(4,2).Range.Text = "71-39="
$t.Cell(4,3).Range.Text = "34+54="
$t.Cell(4,4).Range.Text = "23-20="
$t.Cell(4,5).Range.Text = "59+4="

$t.Cell(5,1).Range.Text = "59-39="
$t.Cell(5,2).Range.Text = "76-13="
$t.Cell(5,3).Range.Text = "40+0="
$t.Cell(5,4).Range.Text = "68-45="
$t.Cell(5,5).Range.Text = "7+26="

$t.Cell(6,1).Range.Text = "1+44="
$t.Cell(6,2).Range.Text = "72-31="
$t.Cell(6,3).Range.Text = "49+17="
$t.Cell(6,4).Range.Text = "81-31="
$t.Cell(6,5).Range.Text = "80+15="

$t.Cell(7,1).Range.Text = "82+3="
$t.Cell(7,2).Range.Text = "75-64="
$t.Cell(7,3).Range.Text = "45-13="
$t.Cell(7,4).Range.Text = "99-81="
$t.Cell(7,5).Range.Text = "57-21="

$t.Cell(8,1).Range.Text = "19-19="
$t.Cell(8,2).Range.Text = "56+1="
$t.Cell(8,3).Range.Text = "3+32="
$t.Cell(8,4).Range.Text = "69-2="
$t.Cell(8,5).Range.Text = "44+18="

$t.Cell(9,1).Range.Text = "50-44="
$t.Cell(9,2).Range.Text = "8+11="
$t.Cell(9,3).Range.Text = "85-59="
$t.Cell(9,4).Range.Text = "69-3="
$t.Cell(9,5).Range.Text = "98-11="

$t.Cell(10,1).Range.Text = "99-19="
$t.Cell(10,2).Range.Text = "66+30="
$t.Cell(10,3).Range.Text = "74-26="
$t.Cell(10,4).Range.Text = "18+40="
$t.Cell(10,5).Range.Text = "9+44="

$t.Cell(11,1).Range.Text = "54-25="
$t.Cell(11,2).Range.Text = "5+80="
$t.Cell(11,3).Range.Text = "15+8="
$t.Cell(11,4).Range.Text = "33-25="
$t.Cell(11,5).Range.Text = "23-11="

$t.Cell(12,1).Range.Text = "20+74="
$t.Cell(12,2).Range.Text = "39-31="
$t.Cell(12,3).Range.Text = "47-14="
$t.Cell(12,4).Range.Text = "6+45="
$t.Cell(12,5).Range.Text = "72+10="

$t.Cell(13,1).Range.Text = "36+35="
$t.Cell(13,2).Range.Text = "58-6="
$t.Cell(13,3).Range.Text = "17+30="
$t.Cell(13,4).Range.Text = "54-33="
$t.Cell(13,5).Range.Text = "62-29="

$t.Cell(14,1).Range.Text = "19+46="
$t.Cell(14,2).Range.Text = "74+6="
$t.Cell(14,3).Range.Text = "62-18="
$t.Cell(14,4).Range.Text = "95-62="
$t.Cell(14,5).Range.Text = "2+49="

$t.Cell(15,1).Range.Text = "94-63="
$t.Cell(15,2).Range.Text = "45+12="
$t.Cell(15,3).Range.Text = "41+17="
$t.Cell(15,4).Range.Text = "33-17="
$t.Cell(15,5).Range.Text = "44-1="

$t.Cell(16,1).Range.Text = "28+1="
$t.Cell(16,2).Range.Text = "34+38="
$t.Cell(16,3).Range.Text = "9+38="
$t.Cell(16,4).Range.Text = "83-82="
$t.Cell(16,5).Range.Text = "34+13="

$t.Cell(17,1).Range.Text = "49+29="
$t.Cell(17,2).Range.Text = "97-70="
$t.Cell(17,3).Range.Text = "23+9="
$t.Cell(17,4).Range.Text = "61+16="
$t.Cell(17,5).Range.Text = "99-17="

$t.Cell(18,1).Range.Text = "15+72="
$t.Cell(18,2).Range.Text = "47-21="
$t.Cell(18,3).Range.Text = "87-80="
$t.Cell(18,4).Range.Text = "23+46="
$t.Cell(18,5).Range.Text = "68-64="

$t.Cell(19,1).Range.Text = "8+25="
$t.Cell(19,2).Range.Text = "26+48="
$t.Cell(19,3).Range.Text = "38-22="
$t.Cell(19,4).Range.Text = "59-14="
$t.Cell(19,5).Range.Text = "48-2="

$t.Cell(20,1).Range.Text = "41+47="
$t.Cell(20,2).Range.Text = "5+20="
$t.Cell(20,3).Range.Text = "46-23="
$t.Cell(20,4).Range.Text = "34+59="
$t.Cell(20,5).Range.Text = "89+7="

